$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.465.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.186.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.96%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.86%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.183.39"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.07%  "
$ws.Range("E9").Value = "  -3.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.141"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.69%  "
$ws.Range("E11").Value = "  -4.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.450"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.11%  "
$ws.Range("E13").Value = "  -6.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.712.85"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.184.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.94%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.508.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "456.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E22").Value = "  -4.77%  "
$ws.Range("E23").Value = "  -4.96%  "
$ws.Range("E24").Value = "  -1.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.87%  "
$ws.Range("E27").Value = "  -3.70%  "
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.93"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.94%  "
$ws.Range("E31").Value = "  -7.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.20%  "
$ws.Range("E33").Value = "  -4.15%  "
$ws.Range("E34").Value = "  -7.20%  "
$ws.Range("E35").Value = "  -5.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.80"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.31"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0694"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0386"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "414.02"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.941.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.13%  "
$ws.Range("B42").Value = "Cosmos"
$ws.Range("C42").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.01"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.32%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.64"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.11%  "
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.112"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.11%  "
$ws.Range("E46").Value = "  -6.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "35.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.67%  "
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.28"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.22%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.43"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.36%  "
$ws.Range("E51").Value = "  -4.48%  "
